# Updates to the Dealer dimension-mapping worksheet: re-point several
# "Source-1" (D:F) rows from the Corporate/Branch source tables over to the
# Products/Dealer source table, clear out now-stale Source-2 (G:I) /
# Source-3 (J:L) mappings, and fill in new Source-2 mappings (Corporate /
# Dealership / Owner_Address columns) where the team found a better match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dealer")

# --- Row 4: DLR_Code now sourced from Products.Dealer instead of Corporate.Dealership ---
$ws.Range("D4").Value = "SEIS732_Team_02_Products"
$ws.Range("E4").Value = "Dealer"
# F4 (DLR_Code) is unchanged
$ws.Range("G4:L4").Value = $null

# --- Row 5: DLR_Name / CD_Name row gains a Source-2 mapping (Corporate.Dealership.ID_Name) ---
$ws.Range("G5").Value = "SEIS732_Team_02_Corporate"
$ws.Range("H5").Value = "Dealership"
$ws.Range("I5").Value = "ID_Name"

# --- Row 6: DLR_Address now sourced from Products.Dealer.DLR_Street instead of Corporate.Branch.BR_Address ---
$ws.Range("D6").Value = "SEIS732_Team_02_Products"
$ws.Range("E6").Value = "Dealer"
$ws.Range("F6").Value = "DLR_Street"
$ws.Range("G6:I6").Value = $null

# --- Row 7: DLR_City now sourced from Products.Dealer.DLR_City instead of Corporate.Branch.BR_City ---
$ws.Range("D7").Value = "SEIS732_Team_02_Products"
$ws.Range("E7").Value = "Dealer"
$ws.Range("F7").Value = "DLR_City"
$ws.Range("G7:I7").Value = $null

# --- Row 8: DLR_County keeps Corporate.Branch.BR_County, gains Source-2 (Corporate.Owner_Address.OADR_County) ---
$ws.Range("G8").Value = "SEIS732_Team_02_Corporate"
$ws.Range("H8").Value = "Owner_Address"
$ws.Range("I8").Value = "OADR_County"

# --- Row 9: DLR_State now sourced from Products.Dealer.DLR_State instead of Corporate.Branch.BR_State ---
$ws.Range("D9").Value = "SEIS732_Team_02_Products"
$ws.Range("E9").Value = "Dealer"
$ws.Range("F9").Value = "DLR_State"
$ws.Range("G9:I9").Value = $null

# --- Row 10: DLR_Country keeps Corporate.Branch.BR_Country, gains Source-2 (Corporate.Owner_Address.OADR_Country) ---
$ws.Range("G10").Value = "SEIS732_Team_02_Corporate"
$ws.Range("H10").Value = "Owner_Address"
$ws.Range("I10").Value = "OADR_Country"

# --- Row 11: DLR_Zip now sourced from Products.Dealer.DLR_Zip instead of Corporate.Branch.BR_Zip ---
$ws.Range("D11").Value = "SEIS732_Team_02_Products"
$ws.Range("E11").Value = "Dealer"
$ws.Range("F11").Value = "DLR_Zip"
$ws.Range("G11:L11").Value = $null

# The new Source-2 values in H/I are noticeably wider than the old ones
# ("Owner_Address" / "OADR_County" / "OADR_Country" vs "Dealer" / "DLR_Code"),
# so re-run AutoFit on those two columns like the author would have.
$ws.Columns("H:H").AutoFit() | Out-Null
$ws.Columns("I:I").AutoFit() | Out-Null

# --- Selection / view cosmetics ---
$ws.Range("G15").Select() | Out-Null
